$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 50 (pushes existing rows 50..158 down to 51..159)
$ws.Rows("50:50").Insert()

# Populate the new row's values
$ws.Range("B50").Value = "strFileHeader29"
$ws.Range("C50").Value = "Field description in exported file"
$ws.Range("D50").Value = "Differentiation algorithm"

# Formatting: left/center aligned, no wrap (matches new cellXfs style for B50)
$ws.Range("B50").HorizontalAlignment = -4131
$ws.Range("B50").VerticalAlignment = -4108
$ws.Range("B50").WrapText = $false

# Formatting: left/center aligned, no wrap, unlocked (matches new cellXfs style for C50:D50)
$ws.Range("C50:D50").HorizontalAlignment = -4131
$ws.Range("C50:D50").VerticalAlignment = -4108
$ws.Range("C50:D50").WrapText = $false
$ws.Range("C50:D50").Locked = $false

# Grow the translation table to cover the new row
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B2:E159"))
